$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.211.13'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '2.647.10'
$ws.Range('E3').Value = '  +2.94%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.62'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.30'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('D9').Value = '2.647.27'
$ws.Range('E9').Value = '  +2.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.106'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.67'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.25%  '
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.41'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('D15').Value = '3.116.87'
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('D16').Value = '63.129.31'
$ws.Range('E16').Value = '  +0.92%  '
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '2.618.93'
$ws.Range('E18').Value = '  +1.77%  '
$ws.Range('E19').Value = '  +1.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '339.42'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.72'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.16%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.92'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.65'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.99%  '
$ws.Range('E26').Value = '  +2.83%  '
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.42'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.64%  '
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '525.51'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +15.42%  '
$ws.Range('E32').Value = '  +12.97%  '
$ws.Range('E33').Value = '  +3.42%  '
$ws.Range('D34').Value = '0.0₃0806'
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '174.50'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.94'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +11.86%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.401'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('E40').Value = '  +7.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '171.32'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +7.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.07'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.74'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.93'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0557'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.631'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.53'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.14%  '
$ws.Range('E51').Value = '  +2.18%  '
